$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A, shifting B:F left to A:E (dimension becomes A1:E3)
$ws.Columns("A:A").Delete()
